$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: insert "mybatis" column before Hibernate, shifting F:I right,
#     and rename the last header from "Mybatis" to "JPA" (net effect of the diff)
$ws.Range("F1").Value = "mybatis"
$ws.Range("G1").Value = "Hibernate"
$ws.Range("H1").Value = "jOOQ"
$ws.Range("I1").Value = "JPA"

# --- Data rows: full new data set (values + reordered/renamed rows + new row 7)
$data = @(
    @("Activiti", "automation", 0, 0.01224739742804654, 0, 0.01000204123290467, 0.0004082465809348847, 0, 0.01285976729944887),
    @("che", "software development", 0.01536643026004728, 0.0007880220646178094, 0, 0.002758077226162333, 0, 0, 0.04137115839243499),
    @("pinpoint", "monitoring", 0, 0.004785863296109952, 0, 0.005276721070069947, 0, 0, 0),
    @("skywalking", "monitoring", 0, 0, 0, 0, 0.0003865481252415926, 0, 0.0003865481252415926),
    @("wildfly", "infrastructure management", 0.002059805383905107, 0.0002130833155763904, 0, 0.01207472121599545, 0.01640741529938206, 0, 0.03821294126003268),
    @("storm", "hpc", 0, 0.0005885815185403178, 0, 0.0008828722778104767, 0, 0.001765744555620953, 0)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $r++
}
